# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets to
# reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10617
$ws1.Range("F3").Value = 296
$ws1.Range("F5").Value = 486
$ws1.Range("F9").Value = 1147
$ws1.Range("F10").Value = 3329
$ws1.Range("F11").Value = 2460
$ws1.Range("F13").Value = 2266
$ws1.Range("F17").Value = 1611
$ws1.Range("F18").Value = 613
$ws1.Range("F20").Value = 274
$ws1.Range("F24").Value = 56
$ws1.Range("F25").Value = 402
$ws1.Range("F26").Value = 17
$ws1.Range("F28").Value = 434
$ws1.Range("F29").Value = 635
$ws1.Range("F30").Value = 61
$ws1.Range("F32").Value = 323
$ws1.Range("F33").Value = 26
$ws1.Range("F34").Value = 1592
$ws1.Range("F35").Value = 660
$ws1.Range("F36").Value = 669
$ws1.Range("F37").Value = 1823
$ws1.Range("F38").Value = 179
$ws1.Range("F39").Value = 479
$ws1.Range("F41").Value = 509
$ws1.Range("F42").Value = 1133
$ws1.Range("F44").Value = 387

# --- 演出 (Performance) sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 21
$ws2.Range("F8").Value = 3

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10617
$ws4.Range("F3").Value = 296
$ws4.Range("F5").Value = 486
$ws4.Range("F8").Value = 1147
$ws4.Range("F9").Value = 3329
$ws4.Range("F10").Value = 2460
$ws4.Range("F11").Value = 2266
$ws4.Range("F14").Value = 1611
$ws4.Range("F15").Value = 613
$ws4.Range("F17").Value = 274
$ws4.Range("F21").Value = 56
$ws4.Range("F22").Value = 402
$ws4.Range("F23").Value = 17
$ws4.Range("F25").Value = 434
$ws4.Range("F26").Value = 635
$ws4.Range("F27").Value = 61
$ws4.Range("F30").Value = 21
$ws4.Range("F32").Value = 323
$ws4.Range("F33").Value = 26
$ws4.Range("F34").Value = 1592
$ws4.Range("F35").Value = 660
$ws4.Range("F37").Value = 669
$ws4.Range("F38").Value = 1823
$ws4.Range("F39").Value = 179
$ws4.Range("F42").Value = 3
$ws4.Range("F43").Value = 479
$ws4.Range("F45").Value = 509
$ws4.Range("F46").Value = 1133
$ws4.Range("F48").Value = 387
